$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.871.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.037.32"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.23"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.58%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("E11").Value = "  +0.76%  "

$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.336.32"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.09"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.762"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.037.40"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.804.69"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.08"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("E24").Value = "  -2.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.28"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.23"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("E28").Value = "  -4.97%  "

$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("E30").Value = "  -6.42%  "

$ws.Range("E31").Value = "  +1.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.45"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.40%  "

$ws.Range("E33").Value = "  +3.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0603"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.49"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.44"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.45%  "

$ws.Range("E38").Value = "  -1.95%  "

$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.539.76"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.01"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0217"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.92"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.50%  "

$ws.Range("E44").Value = "  -1.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0923"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.21%  "

$ws.Range("E46").Value = "  -1.65%  "

$ws.Range("E47").Value = "  -4.76%  "

$ws.Range("E48").Value = "  -1.68%  "

$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.13"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.226.16"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.24%  "
